$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.158.28'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.054.90'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'248.40"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.20%  '
$ws.Range('D6').Value = "'0.665"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.72%  '
$ws.Range('D7').Value = "'57.97"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.22%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.384"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('D10').Value = "'0.0782"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.43%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = "'15.89"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('D13').Value = '2.355.80'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = "'5.73"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = '2.057.32'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = "'18.14"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +17.11%  '
$ws.Range('D18').Value = '37.217.56'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = "'74.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '0.0₃0897'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('D21').Value = "'5.37"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').Value = "'237.71"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('E25').Value = '  -5.95%  '
$ws.Range('D26').Value = "'169.43"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').Value = "'9.40"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('D28').Value = "'20.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').Value = "'0.124"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = "'1.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'4.79"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').Value = "'0.0620"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('D33').Value = "'4.52"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = "'0.0904"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = "'2.28"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').Value = "'1.75"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('E39').Value = '  +12.52%  '
$ws.Range('D40').Value = "'0.103"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.66%  '
$ws.Range('D41').Value = "'5.19"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.60%  '
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'1.14"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = "'17.10"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.11%  '
$ws.Range('D45').Value = "'96.17"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('D46').Value = "'2.46"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('D48').Value = '1.279.48'
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('D49').Value = "'6.85"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('D50').Value = '2.239.06'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = "'43.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.76%  '
